$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" worksheet
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column),
# shifting the old N/O/P columns (Late, heading, Outstanding) one to the right.
$ws.Columns("N").Insert()

# Excel carries the column width of the column to the left (M) into the
# newly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selection to match the resulting state (R10) after the insert.
$null = $ws.Range("R10").Select()
